$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.760.86"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").Value = "1.777.61"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "0.545"
$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "32.12"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "0.289"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").Value = "0.0680"
$ws.Range("E10").Value = "  -5.11%  "

$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").Value = "2.035.40"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = "11.20"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").Value = "1.766.73"
$ws.Range("E14").Value = "  -1.52%  "

$ws.Range("D15").Value = "33.820.71"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "0.610"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("E17").Value = "  -2.58%  "

$ws.Range("D18").Value = "66.58"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").Value = "238.47"
$ws.Range("E19").Value = "  -3.14%  "

$ws.Range("D20").Value = "0.0₃0773"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "10.58"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("D25").Value = "160.50"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").Value = "16.10"
$ws.Range("E26").Value = "  -2.25%  "

$ws.Range("D27").Value = "7.02"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").Value = "0.0511"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "3.59"
$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -2.60%  "

$ws.Range("D35").Value = "1.386.01"

$ws.Range("D36").Value = "0.645"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("E38").Value = "  -0.90%  "

$ws.Range("D39").Value = "2.24"
$ws.Range("E39").Value = "  +5.03%  "

$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "78.22"
$ws.Range("E41").Value = "  -2.44%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.909"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("E43").Value = "  +14.46%  "

$ws.Range("D44").Value = "2.65"
$ws.Range("E44").Value = "  -2.90%  "

$ws.Range("E45").Value = "  +3.31%  "

$ws.Range("D46").Value = "0.0501"
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("E47").Value = "  +11.77%  "

$ws.Range("D48").Value = "107.47"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").Value = "1.933.09"
$ws.Range("E50").Value = "  -0.94%  "
